$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.893.52'
$ws.Range('E2').Value = '  +4.68%  '
$ws.Range('D3').Value = '3.358.38'
$ws.Range('E3').Value = '  +5.03%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '559.58'
$ws.Range('E5').Value = '  +4.06%  '
$ws.Range('D6').Value = '153.41'
$ws.Range('E6').Value = '  +5.92%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('E8').Value = '  +0.87%  '
$ws.Range('D9').Value = '7.53'
$ws.Range('E9').Value = '  +2.57%  '
$ws.Range('E10').Value = '  +4.51%  '
$ws.Range('D11').Value = '0.438'
$ws.Range('E11').Value = '  +1.81%  '
$ws.Range('D12').Value = '3.940.59'
$ws.Range('E12').Value = '  +5.04%  '
$ws.Range('E13').Value = '  +0.23%  '
$ws.Range('D14').Value = '27.18'
$ws.Range('E14').Value = '  +4.35%  '
$ws.Range('D15').Value = '0.0000181'
$ws.Range('E15').Value = '  +3.72%  '
$ws.Range('D16').Value = '62.951.79'
$ws.Range('E16').Value = '  +4.69%  '
$ws.Range('D17').Value = '3.314.85'
$ws.Range('E17').Value = '  +5.31%  '
$ws.Range('E18').Value = '  +4.48%  '
$ws.Range('D19').Value = '13.83'
$ws.Range('E19').Value = '  +5.62%  '
$ws.Range('D20').Value = '8.44'
$ws.Range('E20').Value = '  +1.29%  '
$ws.Range('D21').Value = '390.11'
$ws.Range('E21').Value = '  +1.64%  '
$ws.Range('D22').Value = '0.542'
$ws.Range('E22').Value = '  +2.44%  '
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.09%  '
$ws.Range('E24').Value = '  +0.26%  '
$ws.Range('E25').Value = '  +5.03%  '
$ws.Range('E26').Value = '  +0.61%  '
$ws.Range('E27').Value = '  +7.83%  '
$ws.Range('E28').Value = '  +0.35%  '
$ws.Range('E29').Value = '  +8.21%  '
$ws.Range('E30').Value = '  +4.46%  '
$ws.Range('E31').Value = '  +5.27%  '
$ws.Range('D32').Value = '23.05'
$ws.Range('E32').Value = '  +3.03%  '
$ws.Range('E33').Value = '  +6.94%  '
$ws.Range('E34').Value = '  +2.06%  '
$ws.Range('D35').Value = '1.49'
$ws.Range('E35').Value = '  +9.47%  '
$ws.Range('D36').Value = '160.70'
$ws.Range('E36').Value = '  +2.73%  '
$ws.Range('E37').Value = '  +12.50%  '
$ws.Range('D38').Value = '27.07'
$ws.Range('E38').Value = '  +5.31%  '
$ws.Range('D39').Value = '0.0744'
$ws.Range('E39').Value = '  +4.73%  '
$ws.Range('D40').Value = '2.829.56'
$ws.Range('E40').Value = '  +1.56%  '
$ws.Range('D42').Value = '4.33'
$ws.Range('E42').Value = '  +2.25%  '
$ws.Range('B43').Value = 'Mantle'
$ws.Range('C43').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D43').Value = '0.749'
$ws.Range('E43').Value = '  +3.01%  '
$ws.Range('B44').Value = 'OKB'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D44').Value = '40.77'
$ws.Range('E44').Value = '  +2.46%  '
$ws.Range('E45').Value = '  +5.10%  '
$ws.Range('E46').Value = '  +8.36%  '
$ws.Range('D47').Value = '3.404.64'
$ws.Range('E47').Value = '  +4.96%  '
$ws.Range('E48').Value = '  +2.66%  '
$ws.Range('D49').Value = '6.32'
$ws.Range('E49').Value = '  +2.46%  '
$ws.Range('D50').Value = '0.808'
$ws.Range('E50').Value = '  +1.13%  '
$ws.Range('D51').Value = '283.04'
$ws.Range('E51').Value = '  +6.96%  '
